# SyntheticDataPipeline: only have one Work sameAs and one Person sameAs so all the models don't fall into the same group
#
# This script edits the "CmsWork", "CmsWorkClosing", "CmsWorkOpening" and "CmsPerson"
# sheets so that only a single Work row keeps its wikibase sameAs value (wd:Q937690)
# and only a single Person row keeps its wikibase sameAs value (wd:Q7251). All other
# rows have that cell removed with the remaining cells in the row shifted left by one
# column, and various other data values are regenerated as a consequence (new blank
# node ids, new alternate title / provenance selections, new concept references,
# etc).

$wb = $excel.ActiveWorkbook

# Shifts the contents of a single row left by one column, starting at column
# $deleteCol (which is effectively removed) through $lastCol (which ends up empty).
# Only the given row is touched - other rows in the same columns are left alone.
function Shift-RowLeft {
    param($sheet, $row, $deleteCol, $lastCol)
    for ($c = $deleteCol; $c -lt $lastCol; $c++) {
        $nextVal = $sheet.Cells.Item($row, $c + 1).Value2
        $sheet.Cells.Item($row, $c).Value = $nextVal
    }
    $sheet.Cells.Item($row, $lastCol).ClearContents()
}

# ----------------------------------------------------------------------------------
# CmsWork sheet
# ----------------------------------------------------------------------------------
$wsWork = $wb.Worksheets.Item("CmsWork")

# Row 2 (http://example.com/collection0/work1)
$wsWork.Range("C2").Value = "_:N00fa1d48fadf47d384985bd4c094e573"
$wsWork.Range("D2").Value = "http://example.com/organization4"
$wsWork.Range("E2").Value = "CmsCollection0CmsWork1 alternative title 0"
$wsWork.Range("G2").Value = "CmsCollection0CmsWork1Id0"
Shift-RowLeft $wsWork 2 17 23
$wsWork.Range("S2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"

# Row 3 (http://example.com/collection0/work3)
$wsWork.Range("C3").Value = "_:N722458c0fbe944558e6aa6583a1698a2"
$wsWork.Range("J3").Value = "CmsCollection0CmsWork3 provenance 0"
Shift-RowLeft $wsWork 3 17 23
$wsWork.Range("S3").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:103"

# Row 4 (http://example.com/collection1/work5)
$wsWork.Range("C4").Value = "_:N135efa3017db45f18a70366d925697e0"
$wsWork.Range("D4").Value = "http://example.com/person3"
$wsWork.Range("E4").Value = "CmsCollection1CmsWork5 alternative title 1"
$wsWork.Range("J4").Value = "CmsCollection1CmsWork5 provenance 1"
Shift-RowLeft $wsWork 4 17 23
$wsWork.Range("S4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:106"

# Row 5 (http://example.com/collection1/work7)
$wsWork.Range("C5").Value = "_:N10724d85d5f143d18cda26923e59ba32"
$wsWork.Range("E5").Value = "CmsCollection1CmsWork7 alternative title 1"
$wsWork.Range("G5").Value = "CmsCollection1CmsWork7Id0"
$wsWork.Range("J5").Value = "CmsCollection1CmsWork7 provenance 1"
Shift-RowLeft $wsWork 5 17 23
$wsWork.Range("S5").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:108"

# Row 6 (http://example.com/freestandingwork9)
$wsWork.Range("B6").Value = "_:Na2058460227a45dcb8072e79710e80d2"
$wsWork.Range("C6").Value = "http://example.com/organization1"
$wsWork.Range("I6").Value = "FreestandingWork9 provenance 0"
Shift-RowLeft $wsWork 6 16 22
$wsWork.Range("R6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:109"

# Row 7 (http://example.com/freestandingwork11)
$wsWork.Range("B7").Value = "_:N0ca084e5a9ad4a1090c8d977a58f28db"
$wsWork.Range("C7").Value = "http://example.com/organization4"
$wsWork.Range("D7").Value = "FreestandingWork11 alternative title 1"
$wsWork.Range("F7").Value = "FreestandingWork11Id1"
Shift-RowLeft $wsWork 7 16 22
$wsWork.Range("R7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"

# ----------------------------------------------------------------------------------
# CmsWorkClosing sheet - blank node ids regenerated
# ----------------------------------------------------------------------------------
$wsClosing = $wb.Worksheets.Item("CmsWorkClosing")

$wsClosing.Range("A2").Value = "_:N2061341608bd426f8b6ae3181e9f21e9"
$wsClosing.Range("C2").Value = "_:Nf85a5bbf34694182876717869155ee40"

$wsClosing.Range("A3").Value = "_:Nccfd3f3040d14a1cac2b86658bd611e5"
$wsClosing.Range("C3").Value = "_:Na1e7971f85534bbf9720c059728850c2"

$wsClosing.Range("A4").Value = "_:N90c3265aa22c4710b0505144dfba3aba"
$wsClosing.Range("C4").Value = "_:Nb2d297c21bb84abf83653d98f4eac18a"

$wsClosing.Range("A5").Value = "_:N36abd647d0474e64816c5c5849efc83a"
$wsClosing.Range("C5").Value = "_:N87d2d0c2f94e4a8991fe449606589c66"

$wsClosing.Range("A6").Value = "_:N3219ced4aec54cfcb05605431cbdd8bb"
$wsClosing.Range("C6").Value = "_:N65a6c296b6ef47fbbc47d97d2a750917"

$wsClosing.Range("A7").Value = "_:N1d4c799acde34f7aac0a66af1b7ff43d"
$wsClosing.Range("C7").Value = "_:N8d7f9f3db60a498b965aedd962ad9396"

# ----------------------------------------------------------------------------------
# CmsWorkOpening sheet - blank node ids regenerated (match CmsWorkClosing column C)
# ----------------------------------------------------------------------------------
$wsOpening = $wb.Worksheets.Item("CmsWorkOpening")

$wsOpening.Range("C2").Value = "_:Nf85a5bbf34694182876717869155ee40"
$wsOpening.Range("C3").Value = "_:Na1e7971f85534bbf9720c059728850c2"
$wsOpening.Range("C4").Value = "_:Nb2d297c21bb84abf83653d98f4eac18a"
$wsOpening.Range("C5").Value = "_:N87d2d0c2f94e4a8991fe449606589c66"
$wsOpening.Range("C6").Value = "_:N65a6c296b6ef47fbbc47d97d2a750917"
$wsOpening.Range("C7").Value = "_:N8d7f9f3db60a498b965aedd962ad9396"

# ----------------------------------------------------------------------------------
# CmsPerson sheet - remove the sameAs value from all but the first CmsPerson row
# ----------------------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("CmsPerson")

# Row 2 (http://example.com/person0) keeps its sameAs (wd:Q7251) untouched.

# Row 3 (http://example.com/person1): sameAs is in column E(5), row ends at F(6)
Shift-RowLeft $wsPerson 3 5 6

# Row 4 (http://example.com/person2): sameAs is in column F(6), row ends at G(7)
Shift-RowLeft $wsPerson 4 6 7

# Row 5 (http://example.com/person3): sameAs is in column E(5), row ends at F(6)
Shift-RowLeft $wsPerson 5 5 6

# Row 6 (http://example.com/person4): sameAs is in column F(6), row ends at G(7)
Shift-RowLeft $wsPerson 6 6 7
